$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.881.24"
$ws.Range("E2").Value = "  -5.09%  "

$ws.Range("D3").Value = "2.213.12"
$ws.Range("E3").Value = "  -6.47%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.50"
$ws.Range("E6").Value = "  -8.69%  "

$ws.Range("E7").Value = "  -7.65%  "

$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.560"
$ws.Range("E9").Value = "  -9.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.50"
$ws.Range("E10").Value = "  -10.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.23"
$ws.Range("E11").Value = "  -2.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0826"
$ws.Range("E12").Value = "  -9.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.77"
$ws.Range("E13").Value = "  -8.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  -3.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.861"
$ws.Range("E15").Value = "  -11.74%  "

$ws.Range("D16").Value = "2.554.15"
$ws.Range("E16").Value = "  -6.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.15"
$ws.Range("E17").Value = "  -7.07%  "

$ws.Range("D18").Value = "2.212.30"
$ws.Range("E18").Value = "  -6.22%  "

$ws.Range("D19").Value = "42.767.69"
$ws.Range("E19").Value = "  -5.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.84"
$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("D21").Value = "0.0₃0959"
$ws.Range("E21").Value = "  -9.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.38"
$ws.Range("E22").Value = "  -12.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.15"
$ws.Range("E23").Value = "  -10.89%  "

$ws.Range("E24").Value = "  -9.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.49"
$ws.Range("E25").Value = "  -9.20%  "

$ws.Range("E26").Value = "  -7.84%  "

$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.05"
$ws.Range("E28").Value = "  -9.32%  "

$ws.Range("E29").Value = "  -5.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.29"
$ws.Range("E30").Value = "  -13.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.43"
$ws.Range("E31").Value = "  -8.54%  "

$ws.Range("E32").Value = "  -9.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.81"
$ws.Range("E33").Value = "  -9.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.37"
$ws.Range("E34").Value = "  -7.46%  "

$ws.Range("E35").Value = "  -6.52%  "

$ws.Range("E36").Value = "  +10.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.00"
$ws.Range("E37").Value = "  +14.30%  "

$ws.Range("E38").Value = "  -6.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.44"
$ws.Range("E39").Value = "  -5.88%  "

$ws.Range("E40").Value = "  -11.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.67"
$ws.Range("E41").Value = "  -6.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0325"
$ws.Range("E42").Value = "  -7.78%  "

$ws.Range("D43").Value = "1.892.26"
$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.83"
$ws.Range("E45").Value = "  -9.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.13"
$ws.Range("E46").Value = "  -6.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.206"
$ws.Range("E47").Value = "  -9.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.40"
$ws.Range("E48").Value = "  -3.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "60.44"
$ws.Range("E49").Value = "  -12.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.82"
$ws.Range("E50").Value = "  -8.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.856"
$ws.Range("E51").Value = "  +13.96%  "
